# Adds a new "Italy" worksheet (test data for Italy market), mirroring the
# existing "Portugal" sheet, and updates a couple of pre-existing selections
# that moved as part of this edit (Germany + Slovakia).

$wb = $excel.ActiveWorkbook

# --- Germany: selection moved from B7 to A8:A21 ---------------------------
$wsGermany = $wb.Worksheets.Item("Germany")
$wsGermany.Range("A8:A21").Select()

# --- Slovakia: no longer the active tab; selection becomes "select all" ---
$wsSlovakia = $wb.Worksheets.Item("Slovakia")
$wsSlovakia.Cells.Select()

# --- Italy: new sheet, copied from Portugal (same layout/format) ----------
$wsPortugal = $wb.Worksheets.Item("Portugal")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPortugal.Copy($null, $lastSheet)

$wsItaly = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsItaly.Name = "Italy"

# Fill in the Italy-specific values (NGC reference first, then the market
# name, so new shared-string entries land in the same order as the diff).
$wsItaly.Range("B4").Value = "NGC-3145/T2452/T2455  "
$wsItaly.Range("B2").Value = "Italy Market"

# Portugal's rows 3-4 carry a taller custom height (wrapped text); Italy's
# short values fit on one line, so re-fit them back to the default height.
$wsItaly.Rows.Item(3).AutoFit()
$wsItaly.Rows.Item(4).AutoFit()

# Italy becomes the newly active tab/selection.
$wsItaly.Range("G22").Select()
